$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (shared string) renames
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# Updated GDP (column C) values for each country-year row
$ws.Range("C2").Value = 2870.311589353206
$ws.Range("C3").Value = 697.6889104500298
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 5191.140356354663
$ws.Range("C6").Value = 7772.38875590225
$ws.Range("C7").Value = 1460.056109840828
$ws.Range("C8").Value = 2934.187009790061
$ws.Range("C9").Value = 6923.341809163824
$ws.Range("C10").Value = 7854.952374701078
$ws.Range("C11").Value = 4729.735976516416
$ws.Range("C12").Value = 4547.50930098406
$ws.Range("C13").Value = 1909.084588129339
$ws.Range("C14").Value = 2898.942214704482
$ws.Range("C15").Value = 665.6274194933962
$ws.Range("C16").Value = 1904.346464968814
$ws.Range("C17").Value = 5555.389721901988
$ws.Range("C18").Value = 8082.02845866252
$ws.Range("C19").Value = 1503.870423231357
$ws.Range("C20").Value = 2983.242707849043
$ws.Range("C21").Value = 6967.266654334572
$ws.Range("C22").Value = 8141.91312675461
$ws.Range("C23").Value = 5082.354756663512
$ws.Range("C24").Value = 4633.590358399045
$ws.Range("C25").Value = 1955.461557360978
$ws.Range("C26").Value = 3083.80337578809
$ws.Range("C27").Value = 2965.153206179127
$ws.Range("C28").Value = 1577.487171555845
$ws.Range("C29").Value = 5660.517066940175
$ws.Range("C30").Value = 2024.117324382548
$ws.Range("C31").Value = 4921.848409120176
$ws.Range("C32").Value = 5360.226632400601
$ws.Range("C33").Value = 5122.180090208862
$ws.Range("C34").Value = 3156.723844635973
$ws.Range("C35").Value = 5745.422744292303
$ws.Range("C36").Value = 5642.578115155247
$ws.Range("C37").Value = 2094.024217383061
$ws.Range("C38").Value = 5295.682695961288
$ws.Range("C39").Value = 3212.740625904757
$ws.Range("C40").Value = 5955.175904294275
$ws.Range("C41").Value = 5919.20956823756
$ws.Range("C42").Value = 2201.396847776877
$ws.Range("C43").Value = 951.3148210424945
$ws.Range("C44").Value = 1140.447753778042
$ws.Range("C45").Value = 2286.013198234259
$ws.Range("C46").Value = 720.1523351943922
$ws.Range("C47").Value = 972.7427283025324
$ws.Range("C48").Value = 5412.131646018807
$ws.Range("C49").Value = 3252.634165082374
$ws.Range("C50").Value = 2612.856880840196
$ws.Range("C51").Value = 1627.760281433693
$ws.Range("C52").Value = 3137.260298393558
$ws.Range("C53").Value = 1640.18070024053
$ws.Range("C54").Value = 707.8672001573369
$ws.Range("C55").Value = 711.3043470146426
$ws.Range("C56").Value = 1775.027517189621
$ws.Range("C57").Value = 5996.49696468919
$ws.Range("C58").Value = 6301.696269820412
$ws.Range("C59").Value = 1338.716747746975
$ws.Range("C60").Value = 6103.744960203087
$ws.Range("C61").Value = 1002.388731936373
$ws.Range("C62").Value = 1128.996380299766
$ws.Range("C63").Value = 2361.056581219794
$ws.Range("C64").Value = 726.6520119370772
$ws.Range("C65").Value = 1024.621364522189
$ws.Range("C66").Value = 5330.539154475424
$ws.Range("C67").Value = 3314.741082534716
$ws.Range("C68").Value = 2735.187532014817
$ws.Range("C69").Value = 1625.905825842452
$ws.Range("C70").Value = 3210.869677115934
$ws.Range("C71").Value = 1751.664428859304
$ws.Range("C72").Value = 729.7808175407341
$ws.Range("C73").Value = 731.9993357350996
$ws.Range("C74").Value = 1836.014008604312
$ws.Range("C75").Value = 6114.227214287786
$ws.Range("C76").Value = 6661.86504232374
$ws.Range("C77").Value = 1384.519227335143
$ws.Range("C78").Value = 6249.151036691844
$ws.Range("C79").Value = 1062.040157863007
$ws.Range("C80").Value = 1134.924536209078
$ws.Range("C81").Value = 2425.561644739583
$ws.Range("C82").Value = 747.8284752776283
$ws.Range("C83").Value = 1079.630539001193
$ws.Range("C84").Value = 5176.058803160127
$ws.Range("C85").Value = 3382.563653843273
$ws.Range("C86").Value = 2886.897484630703
$ws.Range("C87").Value = 1644.598009122967
$ws.Range("C88").Value = 3242.636921959078
$ws.Range("C89").Value = 1875.732161108182
$ws.Range("C90").Value = 749.2194349876407
$ws.Range("C91").Value = 1895.214690888655
$ws.Range("C92").Value = 6262.368904654469

# Colony flag corrections (column AL)
$ws.Range("AL3").Value = 1
$ws.Range("AL5").Value = 1
$ws.Range("AL15").Value = 1
$ws.Range("AL17").Value = 1
$ws.Range("AL29").Value = 1
$ws.Range("AL35").Value = 1
$ws.Range("AL40").Value = 1
$ws.Range("AL58").Value = 1
$ws.Range("AL76").Value = 1
